$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.041.80"
$ws.Cells.Item(2, 5).Value = "  -0.49%  "
$ws.Cells.Item(3, 4).Value = "1.828.16"
$ws.Cells.Item(3, 5).Value = "  -0.11%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.008"
$ws.Cells.Item(4, 5).Value = "  -0.21%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "312.62"
$ws.Cells.Item(5, 5).Value = "  -0.22%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.007"
$ws.Cells.Item(6, 5).Value = "  -0.19%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4604"
$ws.Cells.Item(7, 5).Value = "  -1.91%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3696"
$ws.Cells.Item(8, 5).Value = "  +0.84%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07345"
$ws.Cells.Item(9, 5).Value = "  -0.83%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.8705"
$ws.Cells.Item(10, 5).Value = "  -1.16%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07946"
$ws.Cells.Item(11, 5).Value = "  +3.59%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "19.78"
$ws.Cells.Item(12, 5).Value = "  -2.82%  "
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.819.87"
$ws.Cells.Item(13, 5).Value = "  -1.82%  "
$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.338"
$ws.Cells.Item(14, 5).Value = "  -0.94%  "
$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "6.522"
$ws.Cells.Item(15, 5).Value = "  -0.24%  "
$ws.Cells.Item(16, 2).Value = "Litecoin"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "91.67"
$ws.Cells.Item(16, 5).Value = "  -1.60%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.009"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000008875"
$ws.Cells.Item(18, 5).Value = "  +1.64%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "1.005"
$ws.Cells.Item(19, 5).Value = "  -0.44%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.67"
$ws.Cells.Item(21, 4).Value = "26.696.78"
$ws.Cells.Item(21, 5).Value = "  -3.19%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.118"
$ws.Cells.Item(22, 5).Value = "  -2.51%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.57"
$ws.Cells.Item(23, 5).Value = "  -0.59%  "
$ws.Cells.Item(24, 4).Value = "1.927.86"
$ws.Cells.Item(24, 5).Value = "  -7.72%  "
$ws.Cells.Item(25, 5).Value = "  +0.70%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.850"
$ws.Cells.Item(26, 5).Value = "  -1.60%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.54"
$ws.Cells.Item(27, 5).Value = "  +0.02%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.069"
$ws.Cells.Item(28, 5).Value = "  -2.27%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.095"
$ws.Cells.Item(29, 5).Value = "  -1.80%  "
$ws.Cells.Item(30, 5).Value = "  -1.20%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.08872"
$ws.Cells.Item(31, 5).Value = "  -0.70%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.977"
$ws.Cells.Item(32, 5).Value = "  +0.49%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.7307"
$ws.Cells.Item(33, 5).Value = "  -2.00%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.439"
$ws.Cells.Item(34, 5).Value = "  -1.67%  "
$ws.Cells.Item(35, 5).Value = "  -2.61%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.459"
$ws.Cells.Item(36, 5).Value = "  -3.11%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.073"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.05241"
$ws.Cells.Item(38, 5).Value = "  -1.11%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01942"
$ws.Cells.Item(39, 5).Value = "  +0.21%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.941"
$ws.Cells.Item(40, 5).Value = "  +0.03%  "
$ws.Cells.Item(41, 5).Value = "  -2.43%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.5161"
$ws.Cells.Item(42, 5).Value = "  -2.19%  "
$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.1630"
$ws.Cells.Item(43, 5).Value = "  -0.80%  "
$ws.Cells.Item(44, 2).Value = "Frax"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.8604"
$ws.Cells.Item(44, 5).Value = "  -14.84%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "8.208"
$ws.Cells.Item(45, 5).Value = "  -2.27%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.4820"
$ws.Cells.Item(46, 5).Value = "  -1.83%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.008"
$ws.Cells.Item(47, 5).Value = "  -0.19%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "10.20"
$ws.Cells.Item(48, 5).Value = "  -2.42%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "102.09"
$ws.Cells.Item(49, 5).Value = "  -2.23%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.623"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.06228"
$ws.Cells.Item(51, 5).Value = "  -0.79%  "
